$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.6851851851851852
$ws.Range("C2").Value = 0.7551020408163265
$ws.Range("D2").Value = 0.7184466019417475
$ws.Range("E2").Value = 49

# Row 3
$ws.Range("B3").Value = 0.6923076923076923
$ws.Range("C3").Value = 0.6136363636363636
$ws.Range("D3").Value = 0.6506024096385543
$ws.Range("E3").Value = 44

# Row 4
$ws.Range("B4").Value = 0.6881720430107527
$ws.Range("C4").Value = 0.6881720430107527
$ws.Range("D4").Value = 0.6881720430107527
$ws.Range("E4").Value = 0.6881720430107527

# Row 5
$ws.Range("B5").Value = 0.6887464387464388
$ws.Range("C5").Value = 0.684369202226345
$ws.Range("D5").Value = 0.6845245057901509

# Row 6
$ws.Range("B6").Value = 0.68855497350121
$ws.Range("C6").Value = 0.6881720430107527
$ws.Range("D6").Value = 0.6863482744004518

# Row 7
$ws.Range("B7").Value = 0.7631578947368421
$ws.Range("C7").Value = 0.5918367346938775
$ws.Range("D7").Value = 0.6666666666666667
$ws.Range("E7").Value = 49

# Row 8
$ws.Range("B8").Value = 0.6363636363636364
$ws.Range("C8").Value = 0.7954545454545454
$ws.Range("D8").Value = 0.7070707070707071
$ws.Range("E8").Value = 44

# Row 9
$ws.Range("B9").Value = 0.6881720430107527
$ws.Range("C9").Value = 0.6881720430107527
$ws.Range("D9").Value = 0.6881720430107527
$ws.Range("E9").Value = 0.6881720430107527

# Row 10
$ws.Range("B10").Value = 0.6997607655502392
$ws.Range("C10").Value = 0.6936456400742115
$ws.Range("D10").Value = 0.6868686868686869

# Row 11
$ws.Range("B11").Value = 0.7031692133559705
$ws.Range("C11").Value = 0.6881720430107527
$ws.Range("D11").Value = 0.6857825567502988

# Row 12
$ws.Range("B12").Value = 0.6885245901639344
$ws.Range("C12").Value = 0.8571428571428571
$ws.Range("D12").Value = 0.7636363636363637
$ws.Range("E12").Value = 49

# Row 13
$ws.Range("B13").Value = 0.78125
$ws.Range("C13").Value = 0.5681818181818182
$ws.Range("D13").Value = 0.6578947368421052
$ws.Range("E13").Value = 44

# Row 14
$ws.Range("B14").Value = 0.7204301075268817
$ws.Range("C14").Value = 0.7204301075268817
$ws.Range("D14").Value = 0.7204301075268817
$ws.Range("E14").Value = 0.7204301075268817

# Row 15
$ws.Range("B15").Value = 0.7348872950819672
$ws.Range("C15").Value = 0.7126623376623377
$ws.Range("D15").Value = 0.7107655502392345

# Row 16
$ws.Range("B16").Value = 0.732394676537987
$ws.Range("C16").Value = 0.7204301075268817
$ws.Range("D16").Value = 0.7136080670885425

# Row 17
$ws.Range("B17").Value = 0.7857142857142857
$ws.Range("C17").Value = 0.673469387755102
$ws.Range("D17").Value = 0.7252747252747253
$ws.Range("E17").Value = 49

# Row 18
$ws.Range("B18").Value = 0.6862745098039216
$ws.Range("C18").Value = 0.7954545454545454
$ws.Range("D18").Value = 0.736842105263158
$ws.Range("E18").Value = 44

# Row 19
$ws.Range("B19").Value = 0.7311827956989247
$ws.Range("C19").Value = 0.7311827956989247
$ws.Range("D19").Value = 0.7311827956989247
$ws.Range("E19").Value = 0.7311827956989247

# Row 20
$ws.Range("B20").Value = 0.7359943977591037
$ws.Range("C20").Value = 0.7344619666048238
$ws.Range("D20").Value = 0.7310584152689417

# Row 21
$ws.Range("B21").Value = 0.7386675100147585
$ws.Range("C21").Value = 0.7311827956989247
$ws.Range("D21").Value = 0.7307474641939838

# Row 22
$ws.Range("B22").Value = 0.7105263157894737
$ws.Range("C22").Value = 0.5510204081632653
$ws.Range("D22").Value = 0.6206896551724138
$ws.Range("E22").Value = 49

# Row 23
$ws.Range("B23").Value = 0.6
$ws.Range("C23").Value = 0.75
$ws.Range("D23").Value = 0.6666666666666665
$ws.Range("E23").Value = 44

# Row 24
$ws.Range("B24").Value = 0.6451612903225806
$ws.Range("C24").Value = 0.6451612903225806
$ws.Range("D24").Value = 0.6451612903225806
$ws.Range("E24").Value = 0.6451612903225806

# Row 25
$ws.Range("B25").Value = 0.6552631578947368
$ws.Range("C25").Value = 0.6505102040816326
$ws.Range("D25").Value = 0.6436781609195401

# Row 26
$ws.Range("B26").Value = 0.6582342954159592
$ws.Range("C26").Value = 0.6451612903225806
$ws.Range("D26").Value = 0.6424422197503399
